$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "; A FAVOR DE SU +++= documento.menores[0].tratamiento+++ +++= documento.menores[0].nombre+++.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "; A FAVOR DE +++FOR m IN documento.menores+++ +++= `$m.tratamiento+++ +++= `$m.nombre+++ +++END-FOR m.",
    2)

$d.Content.Find.Execute(
    "; A FAVOR DE SU +++= documento.menores[0].tratamiento+++ +++= documento.menores[0].nombre+++, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "; A FAVOR DE SU +++FOR mn IN documento.menores+++ +++= `$mn.tratamiento+++ +++= `$m.nombre+++ +++END-FOR mn+++ ",
    2)
